$wb = $excel.ActiveWorkbook

# --- Sheet "P_valores" ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.124225337243856
$wsP.Range("D2").Value = 0.15309147801399
$wsP.Range("E2").Value = 0.8770976515231426
$wsP.Range("F2").Value = 0.2097823744357821

$wsP.Range("B3").Value = 0.124225337243856
$wsP.Range("D3").Value = 0.9868126371468151
$wsP.Range("E3").Value = 0.2289024156088926
$wsP.Range("F3").Value = 0.6394415397233839

$wsP.Range("B4").Value = 0.15309147801399
$wsP.Range("C4").Value = 0.9868126371468151
$wsP.Range("E4").Value = 0.2866881326726913
$wsP.Range("F4").Value = 0.6780132269713963

$wsP.Range("B5").Value = 0.8770976515231426
$wsP.Range("C5").Value = 0.2289024156088926
$wsP.Range("D5").Value = 0.2866881326726913
$wsP.Range("F5").Value = 0.06144005587143186

$wsP.Range("B6").Value = 0.2097823744357821
$wsP.Range("C6").Value = 0.6394415397233839
$wsP.Range("D6").Value = 0.6780132269713963
$wsP.Range("E6").Value = 0.06144005587143186

# --- Sheet "Estadisticos_DM" ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = -1.635474790239025
$wsE.Range("D2").Value = -1.510746183936936
$wsE.Range("E2").Value = -0.157502792537778
$wsE.Range("F2").Value = -1.314573876725897

$wsE.Range("B3").Value = 1.635474790239025
$wsE.Range("D3").Value = -0.01682629717159442
$wsE.Range("E3").Value = 1.258191216948428
$wsE.Range("F3").Value = 0.4788449270112654

$wsE.Range("B4").Value = 1.510746183936936
$wsE.Range("C4").Value = 0.01682629717159442
$wsE.Range("E4").Value = 1.107632785577087
$wsE.Range("F4").Value = 0.4239932605635565

$wsE.Range("B5").Value = 0.157502792537778
$wsE.Range("C5").Value = -1.258191216948428
$wsE.Range("D5").Value = -1.107632785577087
$wsE.Range("F5").Value = -2.033237266158273

$wsE.Range("B6").Value = 1.314573876725897
$wsE.Range("C6").Value = -0.4788449270112654
$wsE.Range("D6").Value = -0.4239932605635565
$wsE.Range("E6").Value = 2.033237266158273
